$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.676.84'
$ws.Range("E2").Value = '  +1.41%  '

$ws.Range("D3").Value = '1.654.56'
$ws.Range("E3").Value = '  +1.31%  '

$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("E5").Value = '  -0.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.09'
$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3837'
$ws.Range("E7").Value = '  +0.45%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3605'
$ws.Range("E8").Value = '  +0.59%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.08'
$ws.Range("E9").Value = '  -1.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08192'
$ws.Range("E10").Value = '  +0.39%  '

$ws.Range("E11").Value = '  +0.91%  '

$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.40'
$ws.Range("E13").Value = '  +0.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.455'
$ws.Range("E14").Value = '  +0.68%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.441'
$ws.Range("E15").Value = '  +1.96%  '

$ws.Range("E16").Value = '  -0.79%  '

$ws.Range("D17").Value = '1.654.61'
$ws.Range("E17").Value = '  +1.78%  '

$ws.Range("E18").Value = '  +2.72%  '

$ws.Range("E19").Value = '  +1.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.771'
$ws.Range("E20").Value = '  +3.14%  '

$ws.Range("E21").Value = '  +1.35%  '

$ws.Range("E22").Value = '  -0.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.68'
$ws.Range("E23").Value = '  +1.85%  '

$ws.Range("D24").Value = '23.689.29'
$ws.Range("E24").Value = '  +1.46%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.485'
$ws.Range("E25").Value = '  -1.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.023'
$ws.Range("E26").Value = '  -1.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.25'
$ws.Range("E27").Value = '  +0.54%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.25'
$ws.Range("E28").Value = '  +0.94%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.236'
$ws.Range("E29").Value = '  -0.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.13'
$ws.Range("E30").Value = '  +0.29%  '

$ws.Range("D31").Value = '1.839.69'
$ws.Range("E31").Value = '  +1.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.088'
$ws.Range("E32").Value = '  +8.98%  '

$ws.Range("E33").Value = '  +4.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.10'
$ws.Range("E34").Value = '  +5.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.058'
$ws.Range("E35").Value = '  -2.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02805'
$ws.Range("E36").Value = '  +1.88%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2507'
$ws.Range("E37").Value = '  +0.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08804'
$ws.Range("E38").Value = '  +0.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.075'
$ws.Range("E39").Value = '  +2.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06986'
$ws.Range("E40").Value = '  -0.49%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.02'
$ws.Range("E41").Value = '  +6.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6985'
$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.338'
$ws.Range("E43").Value = '  -0.47%  '

$ws.Range("E44").Value = '  +2.78%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6501'
$ws.Range("E45").Value = '  +0.69%  '

$ws.Range("E46").Value = '  -0.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.300'
$ws.Range("E47").Value = '  +1.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.965'
$ws.Range("E48").Value = '  +0.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07888'
$ws.Range("E49").Value = '  -0.69%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.07'
$ws.Range("E50").Value = '  -0.51%  '

$ws.Range("E51").Value = '  -0.74%  '
